$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert E-column total-parcelas for rows 2-10,12,13 to numeric type ---
$ws.Range("E2").Value = 33
$ws.Range("E3").Value = 44
$ws.Range("E4").Value = 39
$ws.Range("E5").Value = 26
$ws.Range("E6").Value = 37
$ws.Range("E7").Value = 45
$ws.Range("E8").Value = 32
$ws.Range("E9").Value = 45
$ws.Range("E10").Value = 34
$ws.Range("E12").Value = 17
$ws.Range("E13").Value = 39

# --- D11 entrada value update ---
$ws.Range("D11").Value = "30190,00"

# --- Row 14 (VC1013) updates ---
$ws.Range("C14").Value = "78000,00"
$ws.Range("D14").Value = "49800,00"
$ws.Range("H14").Value = "1 x R`$ 1645.00`n41 x R`$ 1061.30"

# --- Row 15 (VC1014) updates ---
$ws.Range("C15").Value = "82000,00"
$ws.Range("D15").Value = "50000,00"
$ws.Range("H15").Value = "1 x R`$ 1480.43`n54 x R`$ 933.90"

# --- Row 16 (VC1015) updates ---
$ws.Range("C16").Value = "84600,00"
$ws.Range("D16").Value = "57130,00"
$ws.Range("H16").Value = "46 x R`$ 915.70"

# --- Row 17 (VC1016) updates ---
$ws.Range("C17").Value = "89000,00"
$ws.Range("D17").Value = "46450,00"
$ws.Range("H17").Value = "20 x R`$ 3184.65"

# --- Row 18 (VC1017) updates ---
$ws.Range("B18").Value = "Veículos"
$ws.Range("C18").Value = "92150,00"
$ws.Range("D18").Value = "52107,50"
$ws.Range("H18").Value = "67 x R`$ 1015.55"

# --- Row 19 (VC1018) updates ---
$ws.Range("B19").Value = "Veículos"
$ws.Range("C19").Value = "93300,00"
$ws.Range("D19").Value = "67565,00"
$ws.Range("F19").Value = "Bradesco"
$ws.Range("H19").Value = "1 x R`$ 1210.00`n67 x R`$ 725.62"

# --- Row 20 (VC1019) updates ---
$ws.Range("B20").Value = "Veículos"
$ws.Range("C20").Value = "103400,00"
$ws.Range("D20").Value = "62070,00"
$ws.Range("F20").Value = "Bradesco"
$ws.Range("H20").Value = "1 x R`$ 1860.00`n60 x R`$ 1225.30"

# --- Row 21 (VC1020) updates ---
$ws.Range("B21").Value = "Veículos"
$ws.Range("C21").Value = "116500,00"
$ws.Range("D21").Value = "57825,00"
$ws.Range("F21").Value = "Bradesco"
$ws.Range("H21").Value = "35 x R`$ 2729.07"

# --- New rows 22-30: copy formatting from row 13, then populate ---
$ws.Range("A13:J13").Copy()
$ws.Range("A22:J30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 22 (VC1021)
$ws.Range("A22").Value = "VC1021"
$ws.Range("B22").Value = "Veículos"
$ws.Range("C22").Value = "116500,00"
$ws.Range("D22").Value = "55725,00"
$ws.Range("F22").Value = "Bradesco"
$ws.Range("G22").Value = "Disponível"
$ws.Range("H22").Value = "30 x R`$ 3193.82"

# Row 23 (VC1022)
$ws.Range("A23").Value = "VC1022"
$ws.Range("B23").Value = "Veículos"
$ws.Range("C23").Value = "175400,00"
$ws.Range("D23").Value = "101770,00"
$ws.Range("F23").Value = "Bradesco"
$ws.Range("G23").Value = "Disponível"
$ws.Range("H23").Value = "1 x R`$ 3117.00`n56 x R`$ 2042.00"

# Row 24 (VC1023)
$ws.Range("A24").Value = "VC1023"
$ws.Range("B24").Value = "Veículos"
$ws.Range("C24").Value = "175400,00"
$ws.Range("D24").Value = "113770,00"
$ws.Range("F24").Value = "Bradesco"
$ws.Range("G24").Value = "Disponível"
$ws.Range("H24").Value = "1 x R`$ 2805.00`n48 x R`$ 1707.25"

# Row 25 (VC1024)
$ws.Range("A25").Value = "VC1024"
$ws.Range("B25").Value = "Veículos"
$ws.Range("C25").Value = "351000,00"
$ws.Range("D25").Value = "212550,00"
$ws.Range("F25").Value = "Bradesco"
$ws.Range("G25").Value = "Disponível"
$ws.Range("H25").Value = "1 x R`$ 5922.00`n52 x R`$ 3750.00"

# Row 26 (VC1025)
$ws.Range("A26").Value = "VC1025"
$ws.Range("B26").Value = "Veículos"
$ws.Range("C26").Value = "526500,00"
$ws.Range("D26").Value = "316325,00"
$ws.Range("F26").Value = "Bradesco"
$ws.Range("G26").Value = "Disponível"
$ws.Range("H26").Value = "1 x R`$ 8319.00`n59 x R`$ 5214.00"

# Row 27 (VC1026)
$ws.Range("A27").Value = "VC1026"
$ws.Range("B27").Value = "Imóveis"
$ws.Range("C27").Value = "99600,00"
$ws.Range("D27").Value = "54880,00"
$ws.Range("F27").Value = "Bradesco"
$ws.Range("G27").Value = "Disponível"
$ws.Range("H27").Value = "100 x R`$ 616.00"

# Row 28 (VC1027)
$ws.Range("A28").Value = "VC1027"
$ws.Range("B28").Value = "Imóveis"
$ws.Range("C28").Value = "140000,00"
$ws.Range("D28").Value = "84000,00"
$ws.Range("F28").Value = "Porto Seguro"
$ws.Range("G28").Value = "Disponível"
$ws.Range("H28").Value = "60 x R`$ 1356.82"

# Row 29 (VC1028)
$ws.Range("A29").Value = "VC1028"
$ws.Range("B29").Value = "Imóveis"
$ws.Range("C29").Value = "404000,00"
$ws.Range("D29").Value = "247200,00"
$ws.Range("F29").Value = "Porto Seguro"
$ws.Range("G29").Value = "Disponível"
$ws.Range("H29").Value = "158 x R`$ 1898.72"

# Row 30 (VC1029)
$ws.Range("A30").Value = "VC1029"
$ws.Range("B30").Value = "Imóveis"
$ws.Range("C30").Value = "544000,00"
$ws.Range("D30").Value = "329200,00"
$ws.Range("F30").Value = "Porto Seguro"
$ws.Range("G30").Value = "Disponível"
$ws.Range("H30").Value = "60 x R`$ 3255.00`n98 x R`$ 1898.00"

# --- Force text type for pure-numeric-looking E-column values via scratch cell round-trip ---
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"

$textCells = @(
    @{ Ref = "E14"; Text = "42" },
    @{ Ref = "E15"; Text = "55" },
    @{ Ref = "E16"; Text = "46" },
    @{ Ref = "E17"; Text = "20" },
    @{ Ref = "E18"; Text = "67" },
    @{ Ref = "E19"; Text = "68" },
    @{ Ref = "E20"; Text = "61" },
    @{ Ref = "E21"; Text = "35" },
    @{ Ref = "E22"; Text = "30" },
    @{ Ref = "E23"; Text = "57" },
    @{ Ref = "E24"; Text = "49" },
    @{ Ref = "E25"; Text = "53" },
    @{ Ref = "E26"; Text = "60" },
    @{ Ref = "E27"; Text = "100" },
    @{ Ref = "E28"; Text = "60" },
    @{ Ref = "E29"; Text = "158" },
    @{ Ref = "E30"; Text = "158" }
)

foreach ($item in $textCells) {
    $scratch.Value = $item.Text
    $scratch.Copy()
    $ws.Paste($ws.Range($item.Ref))
}

# --- Cleanup scratch column ---
$ws.Columns("ZZ").Delete() | Out-Null
$excel.CutCopyMode = $false
